$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pH readings for "week 03" (col D) and "week 04" (col E) for rows 6-9
$ws.Range("D6").Value = 6.36
$ws.Range("E6").Value = 6.61

$ws.Range("D7").Value = 6.79
$ws.Range("E7").Value = 6.79

$ws.Range("D8").Value = 7.01
$ws.Range("E8").Value = 6.74

$ws.Range("D9").Value = 6.21
$ws.Range("E9").Value = 6.46

# The newly entered "week 03" values picked up an explicit (General) alignment
# touch while typing them in, which stamps its own cell style record
$ws.Range("D6:D9").HorizontalAlignment = 1

# Remove the old Notes row (row 12) entirely, shifting the used range back up
$ws.Rows(12).Delete()

# Move the active selection to C10 (matches where the author left off editing)
$ws.Range("C10").Select()
